# Aston Villa_stats.xlsx edit script
# 1) Rename the per-stat-category sheet tabs to their spaced / punctuated
#    display names.
# 2) Bump every player's "age-days since birthday" value (column E,
#    format "YY-DDD") in the stat sheets by one day, reflecting the
#    workbook being regenerated a day later.

$wb = $excel.ActiveWorkbook

# --- 1) Rename sheet tabs -------------------------------------------------
$renames = @{
    "StandardStats"    = "Standard Stats"
    "ShootingStats"     = "Shooting Stats"
    "PassingStats"      = "Passing Stats"
    "PassTypes"         = "Pass Types"
    "GoalShotCreation"  = "Goal & Shot Creation"
    "DefensiveActions"  = "Defensive Actions"
    "PlayingTime"       = "Playing Time"
    "MiscStats"         = "Miscellaneous Stats"
}

foreach ($oldName in $renames.Keys) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $renames[$oldName]
}

# --- 2) Increment the "Age" column (E) day component by one --------------
# Values look like "27-317" (27 years, 317 days) -> "27-318".
# Applies to the player rows (4..N) on every stats sheet; "Standard Stats"
# and "Playing Time" list a few extra unused-player rows (up to row 37),
# the rest stop at row 31.

$statSheets = @(
    "Standard Stats",
    "Shooting Stats",
    "Passing Stats",
    "Pass Types",
    "Goal & Shot Creation",
    "Defensive Actions",
    "Possession",
    "Playing Time",
    "Miscellaneous Stats"
)

foreach ($name in $statSheets) {
    $ws = $wb.Worksheets.Item($name)

    $lastRow = 31
    if ($name -eq "Standard Stats" -or $name -eq "Playing Time") {
        $lastRow = 37
    }

    for ($r = 4; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 5)
        $val = $cell.Value()
        if ($val -match "^(\d+)-(\d+)$") {
            $years = $matches[1]
            $days = [int]$matches[2] + 1
            $newVal = "{0}-{1:D3}" -f $years, $days
            $cell.Value = $newVal
        }
    }
}
